# Answer the 3rd part of the 1st question ("What are potential solutions?").
#
# Before:
#   "...A constraint is the size..."              <- carries the _GoBack bookmark at its end
#   "What are potential solutions?" (bold)
#   ""                                             <- empty paragraph right after it
#   "Evaluate each solution?" (bold)
#
# After:
#   "...A constraint is the size..."              <- bookmark removed
#   "What are potential solutions?" (bold)         <- unchanged
#   "Potential solutions are to get a bigger boat, to find additional help, or to keep
#    the cat separate from the parrot and the parrot separate from the seeds."
#                                                   <- new text, now carries the _GoBack
#                                                      bookmark at its end
#   "Evaluate each solution?" (bold)               <- unchanged

$d = $word.ActiveDocument

# 1. Remove the existing "_GoBack" bookmark from the end of the "A constraint..." paragraph
#    (it will be re-created after the new answer text below).
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# 2. Find the "What are potential solutions?" paragraph; the paragraph right after it is the
#    (currently empty) one that should receive the new answer.
$searchRange = $d.Content
$null = $searchRange.Find.Execute("What are potential solutions?", $false, $false, $false,
                                   $false, $false, $true, 1, $false, "", 0)
$solutionsPara = $searchRange.Paragraphs(1)
$answerParagraph = $solutionsPara.Next()

# 3. Fill that empty paragraph with the new answer text. A sentinel character is appended
#    temporarily so the bookmark re-inserted below can be anchored precisely at "end of the
#    real text" - a zero-length range placed exactly before the paragraph mark snaps to the
#    paragraph boundary instead of the desired offset.
$answerText = "Potential solutions are to get a bigger boat, to find additional help, or to keep the cat separate from the parrot and the parrot separate from the seeds."
$sentinel = "#"
$answerParagraph.Range.Text = $answerText + $sentinel

# 4. Re-create the "_GoBack" bookmark right after the real text (i.e. right before the
#    sentinel character).
$answerParagraph = $solutionsPara.Next()
$bmRange = $answerParagraph.Range.Duplicate
$null = $bmRange.MoveEnd(1, -2)
$bmRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $bmRange)

# 5. Strip the sentinel character back out again.
$answerParagraph = $solutionsPara.Next()
$sentinelRange = $answerParagraph.Range.Duplicate
$null = $sentinelRange.MoveEnd(1, -1)
$null = $sentinelRange.MoveStart(1, $sentinelRange.End - $sentinelRange.Start - 1)
$sentinelRange.Text = ""

Write-Output "Answer paragraph now reads: $($solutionsPara.Next().Range.Text)"
